$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the season record columns, matching the
# existing header row formatting (bold, bordered, centered).
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 29).Value = 68
    $ws.Cells.Item($r, 30).Value = 94
    $ws.Cells.Item($r, 31).Value = 0
}
